# "added prob 12 in common part"
# Adds a new "common part" problem (d0013-d0017) as rows 53-57 on Sheet1,
# plus five new "key" cells (C21:C25) that pair up with existing building
# blocks in rows 21-25, and shrinks column C's width / updates the saved
# view position to show the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New "key" values for existing rows 21-25 (column C was empty there) ---
$ws.Range("C21").Value = "두 지수함수와 기울기 2인 직선;"
$ws.Range("C22").Value = "`$\sqrt{1+m^{2}} | x_{1} - x_{2} |`$"
$ws.Range("C23").Value = "곡선 위의 점이면서 직선 위의 점;"
$ws.Range("C24").Value = "미정계수 a, b, c, d;"
$ws.Range("C25").Value = "(0, 0); (1, 2);"

# --- New rows 53-57: prob 12 building blocks (common part) ---
$ws.Range("A53").Value = "d0013"
$ws.Range("B53").Value = "우변0인 항등식의 좌변을 인수분해합니다."

$ws.Range("A54").Value = "d0014"
$ws.Range("B54").Value = "인수분해된 식을 이용해 가능한 함수의 식을 찾아 줍니다."
$ws.Range("C54").Value = "`$f(x)`$에 대한 삼차식 인수분해;"

$ws.Range("A55").Value = "d0015"
$ws.Range("B55").Value = "최대와 최소의 조건에 맞게 실수 전체의 집합에서 연속하도록하는 구간 함수를 구간별로 구해줍니다."
$ws.Range("C55").Value = "최대1 최소0;"

$ws.Range("A56").Value = "d0016"
$ws.Range("B56").Value = "구해진 결정된 함수에서 각각의 함숫값을 구합니다."
$ws.Range("C56").Value = "구간함수의 함숫값;"

$ws.Range("A57").Value = "d0017"
$ws.Range("B57").Value = "각각의 함숫값을 이용해서 문제에서 요구하는 것을 구합니다."
$ws.Range("C57").Value = "함숫값들의 합;"

# --- Column C got narrower (was 100 chars wide) ---
$ws.Columns.Item(3).ColumnWidth = 65.57

# --- Update the saved view: scroll down and select the last new cell ---
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C57").Select()
